$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-11 Sunday" "2024-02-12 Monday"

Replace-Text "22×86=" "22×17="
Replace-Text "97×91=" "44×84="
Replace-Text "93×55=" "77×81="
Replace-Text "44×55=" "20×33="
Replace-Text "84×38=" "56×26="
Replace-Text "40×64=" "69×16="
Replace-Text "85×88=" "44×62="
Replace-Text "11×36=" "86×88="
Replace-Text "13×90=" "48×65="
Replace-Text "52×38=" "26×63="
Replace-Text "50×90=" "44×64="
Replace-Text "90×51=" "24×39="
Replace-Text "98×67=" "19×93="
Replace-Text "87×69=" "40×77="
Replace-Text "15×77=" "17×61="
Replace-Text "24×68=" "70×74="
Replace-Text "92×24=" "17×94="
Replace-Text "94×16=" "63×29="
Replace-Text "24×95=" "61×17="
Replace-Text "71×24=" "96×95="
Replace-Text "22×56=" "52×46="
Replace-Text "89×46=" "75×55="
Replace-Text "99×29=" "88×31="
Replace-Text "54×27=" "84×66="
Replace-Text "94×43=" "50×95="
